# Save the data for each textbox ("day 1".."day 5") on each day.
# Columns A,C,E,G,I,K,M hold the day-number label textboxes; columns
# B,D,F,H,J,L,N hold the paired True/False checkbox state for that day.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "A" = @(" 1", " 2", " 3", "4", " 5")
    "B" = @("True", "True", "False", "False", "False")
    "C" = @("1", "2", "3", "4", "5")
    "D" = @("False", "True", "False", "False", "False")
    "E" = @("1", "2", "3", "4", "5")
    "F" = @("False", "False", "True", "False", "False")
    "G" = @(" 1", " 2", " 3", "4", " 5")
    "H" = @("False", "True", "False", "True", "False")
    "I" = @("1", "2", "3", "4", "5")
    "J" = @("False", "False", "False", "False", "False")
    "K" = @(" 1", " 2", " 3", " 4", " 5")
    "L" = @("False", "False", "False", "False", "False")
    "M" = @(" 1", " 2", " 3", " 4", " 5")
    "N" = @("False", "False", "False", "False", "False")
}

foreach ($col in $data.Keys) {
    $values = $data[$col]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $row = $i + 1
        # Leading apostrophe forces Excel to store the value as literal
        # text (inline/shared string) instead of coercing "True"/"False"
        # to a boolean or numeric-looking strings ("1") to a number.
        $ws.Range("$col$row").Value = "'" + $values[$i]
    }
}

# New bestFit-ish widths for the newly populated label columns.
$ws.Range("C1").ColumnWidth = 15.25
$ws.Range("E1").ColumnWidth = 15.25
$ws.Range("G1").ColumnWidth = 16.45
$ws.Range("I1").ColumnWidth = 17.5

# Selection moves to G10 as recorded in the saved view state.
$ws.Range("G10").Select()
